# Updates betting odds values on Sheet1 for row 7 and row 14
# as described by the source diff (Atualizando o arquivo XLSX).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---- Row 7 updates ----
$ws.Range("H7").Value  = 2.87
$ws.Range("I7").Value  = 3.35
$ws.Range("J7").Value  = 2.9
$ws.Range("L7").Value  = 3.8
$ws.Range("O7").Value  = 1.39
$ws.Range("P7").Value  = 2.55
$ws.Range("T7").Value  = 2.4
$ws.Range("W7").Value  = 6.3
$ws.Range("X7").Value  = 10
$ws.Range("AA7").Value = 21
$ws.Range("AB7").Value = 35
$ws.Range("AC7").Value = 7.3
$ws.Range("AF7").Value = 75
$ws.Range("AG7").Value = 9
$ws.Range("AH7").Value = 18
$ws.Range("AI7").Value = 11.5
$ws.Range("AK7").Value = 32
$ws.Range("AL7").Value = 40
$ws.Range("AO7").Value = 12
$ws.Range("AP7").Value = 21
$ws.Range("AQ7").Value = 50
$ws.Range("AR7").Value = 90
$ws.Range("AS7").Value = 300
$ws.Range("AT7").Value = 2.37
$ws.Range("AU7").Value = 6.8
$ws.Range("AW7").Value = 5.1
$ws.Range("AX7").Value = 18.5
$ws.Range("AY7").Value = 24
$ws.Range("AZ7").Value = 90
$ws.Range("BA7").Value = 120
$ws.Range("BB7").Value = 300

# ---- Row 14 updates ----
$ws.Range("J14").Value  = 3.45
$ws.Range("R14").Value  = 1.6
$ws.Range("S14").Value  = 1.42
$ws.Range("T14").Value  = 2.47
$ws.Range("W14").Value  = 8.75
$ws.Range("X14").Value  = 15.5
$ws.Range("Z14").Value  = 40
$ws.Range("AC14").Value = 7.8
$ws.Range("AD14").Value = 5.7
$ws.Range("AG14").Value = 7
$ws.Range("AN14").Value = 4.85
$ws.Range("AP14").Value = 21
$ws.Range("AT14").Value = 2.45
$ws.Range("AU14").Value = 6.6
